$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23, shifting existing rows 23:61 down to 24:62
$ws.Rows.Item(23).Insert($excel.XlInsertShiftDirection.xlShiftDown)

# Fill in the new row 23 with data (same Mercado/Producto context as surrounding rows)
$ws.Range("A23").Value = 10
$ws.Range("B23").Value = "Vega Modelo de Temuco"
$ws.Range("C23").Value = "La Araucanía"
$ws.Range("D23").Value = 44467
$ws.Range("E23").Value = 9
$ws.Range("F23").Value = "Fruta"
$ws.Range("G23").Value = 100107
$ws.Range("H23").Value = "Otros"
$ws.Range("I23").Value = 100107002
$ws.Range("J23").Value = "Chirimoya"
$ws.Range("K23").Value = "Sin especificar"
$ws.Range("L23").Value = "Primera"
$ws.Range("M23").Value = 30
$ws.Range("N23").Value = 3200
$ws.Range("O23").Value = 3200
$ws.Range("P23").Value = 3200
$ws.Range("Q23").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R23").Value = "Provincia del Elquí"
$ws.Range("S23").Value = 3200
$ws.Range("T23").Value = 1
